$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: 山贼强盗
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "山贼强盗"
$ws.Range("C4").Value = "山贼强盗"
$ws.Range("D4").Value = "山贼强盗"
$ws.Range("E4").Value = 80
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 8
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 2
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 2

# Row 5: 强盗头子
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "强盗头子"
$ws.Range("C5").Value = "强盗头子"
$ws.Range("D5").Value = "强盗头子"
$ws.Range("E5").Value = 200
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 20
$ws.Range("J5").Value = 20
$ws.Range("K5:P5").ClearContents()

# Row 6: 花妖
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "花妖"
$ws.Range("C6").Value = "花妖"
$ws.Range("D6").Value = "花妖"
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 7
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 1
$ws.Range("N6").Value = 2
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 2

# Row 7: 山精
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "山精"
$ws.Range("C7").Value = "山精"
$ws.Range("D7").Value = "山精"
$ws.Range("E7").Value = 50
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 7
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 1
$ws.Range("N7").Value = 2
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = 2

# Row 8: 鬼怪
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "鬼怪"
$ws.Range("C8").Value = "鬼怪"
$ws.Range("D8").Value = "鬼怪"
$ws.Range("E8").Value = 50
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 7
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 1
$ws.Range("N8").Value = 2
$ws.Range("O8").Value = 1
$ws.Range("P8").Value = 2

# Row 9: 食人花
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "食人花"
$ws.Range("C9").Value = "食人花"
$ws.Range("D9").Value = "食人花"
$ws.Range("E9").Value = 150
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 15
$ws.Range("H9").Value = 27
$ws.Range("I9").Value = 14
$ws.Range("J9").Value = 22
$ws.Range("K9:P9").ClearContents()

# Row 10: 白骨士兵
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "白骨士兵"
$ws.Range("C10").Value = "白骨士兵"
$ws.Range("D10").Value = "白骨士兵"
$ws.Range("E10").Value = 70
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 4
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 10
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 2
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 1
$ws.Range("P10").Value = 1

# Row 11: 白骨弓箭手
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "白骨弓箭手"
$ws.Range("C11").Value = "白骨弓箭手"
$ws.Range("D11").Value = "白骨弓箭手"
$ws.Range("E11").Value = 100
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 13
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 3
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 8
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 2
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 1
$ws.Range("P11").Value = 1

# Row 12: 白骨将军
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "白骨将军"
$ws.Range("C12").Value = "白骨将军"
$ws.Range("D12").Value = "白骨将军"
$ws.Range("E12").Value = 300
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 54
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 39
$ws.Range("J12").Value = 30
$ws.Range("K12:P12").ClearContents()

# Row 13: 白骨之王
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "白骨之王"
$ws.Range("C13").Value = "白骨之王"
$ws.Range("D13").Value = "白骨之王"
$ws.Range("E13").Value = 1000
$ws.Range("F13").Value = 200
$ws.Range("G13").Value = 100
$ws.Range("H13").Value = 40
$ws.Range("I13").Value = 95
$ws.Range("J13").Value = 80
$ws.Range("K13:P13").ClearContents()

# Row 14: 海妖
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "海妖"
$ws.Range("C14").Value = "海妖"
$ws.Range("D14").Value = "海妖"
$ws.Range("E14").Value = 30
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 2
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = 9
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 1
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 1
$ws.Range("P14").Value = 1

# Row 15: 水鬼
$ws.Range("A15").Value = 12
$ws.Range("B15").Value = "水鬼"
$ws.Range("C15").Value = "水鬼"
$ws.Range("D15").Value = "水鬼"
$ws.Range("E15").Value = 80
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 1
$ws.Range("K15").Value = 9
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 4
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 1
$ws.Range("P15").Value = 1

# Row 16: 鱼人
$ws.Range("A16").Value = 13
$ws.Range("B16").Value = "鱼人"
$ws.Range("C16").Value = "鱼人"
$ws.Range("D16").Value = "鱼人"
$ws.Range("E16").Value = 60
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 11
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 2
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 1
$ws.Range("P16").Value = 1

# Row 17: 嗜血鱼人
$ws.Range("A17").Value = 14
$ws.Range("B17").Value = "嗜血鱼人"
$ws.Range("C17").Value = "嗜血鱼人"
$ws.Range("D17").Value = "嗜血鱼人"
$ws.Range("E17").Value = 450
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 32
$ws.Range("J17").Value = 40
$ws.Range("K17:P17").ClearContents()

# Row 18: 叛军
$ws.Range("A18").Value = 15
$ws.Range("B18").Value = "叛军"
$ws.Range("C18").Value = "叛军"
$ws.Range("D18").Value = "叛军"
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 33
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 9
$ws.Range("J18").Value = 9
$ws.Range("K18").Value = 12
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 3
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 3
$ws.Range("P18").Value = 3

# Row 19: 叛军首领
$ws.Range("A19").Value = 16
$ws.Range("B19").Value = "叛军首领"
$ws.Range("C19").Value = "叛军首领"
$ws.Range("D19").Value = "叛军首领"
$ws.Range("E19").Value = 2250
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 167
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 95
$ws.Range("J19").Value = 70
$ws.Range("K19:P19").ClearContents()

# Row 20: 幽灵树精
$ws.Range("A20").Value = 17
$ws.Range("B20").Value = "幽灵树精"
$ws.Range("C20").Value = "幽灵树精"
$ws.Range("D20").Value = "幽灵树精"
$ws.Range("E20").Value = 330
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 25
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 18
$ws.Range("J20").Value = 27
$ws.Range("K20").Value = 15
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 2
$ws.Range("N20").Value = 4
$ws.Range("O20").Value = 2
$ws.Range("P20").Value = 3

# Row 21: 幽灵树王
$ws.Range("A21").Value = 18
$ws.Range("B21").Value = "幽灵树王"
$ws.Range("C21").Value = "幽灵树王"
$ws.Range("D21").Value = "幽灵树王"
$ws.Range("E21").Value = 2750
$ws.Range("F21").Value = 1000
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = 114
$ws.Range("I21").Value = 54
$ws.Range("J21").Value = 99
$ws.Range("K21:P21").ClearContents()

# Clear the bold/colored style previously applied to B4:D5 (now plain, unused font removed)
$ws.Range("B4:D5").ClearFormats()

# Update selection to match target view state
$ws.Range("P24").Select()